$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UpdatedLayout")

# ---------------------------------------------------------------
# 1) Remove table blocks that are fully vacated in the new layout
#    (drops both content and the grey "table name" formatting)
# ---------------------------------------------------------------
$ws.Range("B23:D23").Clear()
$ws.Range("B28:D28").Clear()
$ws.Range("B32:D32").Clear()
$ws.Range("B39").Clear()

# ---------------------------------------------------------------
# 2) Drop remaining stray cells with default formatting that have
#    no counterpart in the new layout
# ---------------------------------------------------------------
$ws.Range("E19").ClearContents()
$ws.Range("E24").ClearContents()
$ws.Range("E26").ClearContents()
$ws.Range("E30").ClearContents()
$ws.Range("E33").ClearContents()
$ws.Range("C41").ClearContents()
$ws.Range("D41").ClearContents()
$ws.Range("C44").ClearContents()
$ws.Range("D44").ClearContents()
$ws.Range("E44").ClearContents()
$ws.Range("E45").ClearContents()

# ---------------------------------------------------------------
# 3) Write the updated layout content
# ---------------------------------------------------------------
$ws.Range("B2").Value = "According to Issue#1 - new Layout with foreign Keys."
$ws.Range("B4").Value = "Table"
$ws.Range("C4").Value = "Keys"
$ws.Range("D4").Value = "Type"
$ws.Range("E4").Value = "Comment"
$ws.Range("B5").Value = "Creator"
$ws.Range("C5").Value = "CreatorID"
$ws.Range("D5").Value = "INTEGER"
$ws.Range("E5").Value = "Changed User -> Creator, since we're not talking about users but people who made stuff."
$ws.Range("C6").Value = "Surname"
$ws.Range("D6").Value = "VARCHAR"
$ws.Range("C7").Value = "Name"
$ws.Range("D7").Value = "VARCHAR"
$ws.Range("C8").Value = "Studies"
$ws.Range("D8").Value = "VARCHAR"
$ws.Range("B10").Value = "Course"
$ws.Range("C10").Value = "CourseID"
$ws.Range("D10").Value = "INTEGER"
$ws.Range("C11").Value = "Programme"
$ws.Range("D11").Value = "VARCHAR"
$ws.Range("C12").Value = "Course"
$ws.Range("D12").Value = "VARCHAR"
$ws.Range("B14").Value = "SceneAsset"
$ws.Range("C14").Value = "AssetID"
$ws.Range("D14").Value = "INTEGER"
$ws.Range("C15").Value = "(FK) Creator"
$ws.Range("D15").Value = "INTEGER"
$ws.Range("C16").Value = "(FK) Course"
$ws.Range("D16").Value = "INTEGER"
$ws.Range("C17").Value = "Filename"
$ws.Range("D17").Value = "VARCHAR"
$ws.Range("C18").Value = "Filetype"
$ws.Range("D18").Value = "VARCHAR"
$ws.Range("C19").Value = "Date"
$ws.Range("D19").Value = "DATE"
$ws.Range("C20").Value = "Link"
$ws.Range("D20").Value = "VARCHAR"
$ws.Range("E20").Value = "(might as well be BLOB later on - or not. We'll see)"
$ws.Range("C21").Value = "Thumbnail"
$ws.Range("D21").Value = "VARCHAR"
$ws.Range("E21").Value = "Will create a link to a jpg/png/... file upon submission"
$ws.Range("C22").Value = "Deleted"
$ws.Range("D22").Value = "BOOLEAN"
$ws.Range("E22").Value = "(Asset deleted? -> True ==> Purge Database in regular intervals)"
$ws.Range("B24").Value = "LightAsset"
$ws.Range("C24").Value = "AssetID"
$ws.Range("D24").Value = "INTEGER"
$ws.Range("C25").Value = "Type"
$ws.Range("D25").Value = "VARCHAR"
$ws.Range("E25").Value = "(Spot, Point, Directional, Area...)"
$ws.Range("C26").Value = "Power"
$ws.Range("D26").Value = "FLOAT"
$ws.Range("C27").Value = "Color"
$ws.Range("D27").Value = "VARCHAR"
$ws.Range("E27").Value = "(will be Vec3 as string -> parse as you go)"
$ws.Range("B29").Value = "MarkerAsset"
$ws.Range("C29").Value = "AssetID"
$ws.Range("D29").Value = "INTEGER"
$ws.Range("C30").Value = "Name"
$ws.Range("D30").Value = "VARCHAR"
$ws.Range("C31").Value = "Link"
$ws.Range("D31").Value = "VARCHAR"
$ws.Range("E31").Value = "(Similar to SceneAsset)"
$ws.Range("B33").Value = "Scene"
$ws.Range("C33").Value = "SceneID"
$ws.Range("D33").Value = "INTEGER"
$ws.Range("C34").Value = "Name"
$ws.Range("D34").Value = "VARCHAR"
$ws.Range("B36").Value = "Anchor"
$ws.Range("C36").Value = "AnchorID"
$ws.Range("D36").Value = "INTEGER"
$ws.Range("F36").Value = "*Shared Primary Key and Class Table Inheritance"
$ws.Range("C37").Value = "SceneID"
$ws.Range("D37").Value = "INTEGER"
$ws.Range("C38").Value = "Transform"
$ws.Range("D38").Value = "VARCHAR"
$ws.Range("E38").Value = "(will be Vec3 as string -> parse as you go)"
$ws.Range("C39").Value = "Rotation"
$ws.Range("D39").Value = "VARCHAR"
$ws.Range("E39").Value = "(will be Vec4 as string -> parse as you go)"
$ws.Range("C40").Value = "Scale"
$ws.Range("D40").Value = "VARCHAR"
$ws.Range("E40").Value = "(will be Vec3 as string -> parse as you go)"
$ws.Range("B42").Value = "LightAnchor"
$ws.Range("C42").Value = "AnchorID"
$ws.Range("D42").Value = "INTEGER"
$ws.Range("C43").Value = "AssetID"
$ws.Range("D43").Value = "INTEGER"
$ws.Range("E43").Value = "LightAsset Table"
$ws.Range("B45").Value = "MarkerAnchor"
$ws.Range("C45").Value = "AnchorID"
$ws.Range("D45").Value = "INTEGER"
$ws.Range("C46").Value = "AssetID"
$ws.Range("D46").Value = "INTEGER"
$ws.Range("E46").Value = "MarkerAssetTable"
$ws.Range("B48").Value = "SceneAssetAnchor"
$ws.Range("C48").Value = "AnchorID"
$ws.Range("D48").Value = "INTEGER"
$ws.Range("C49").Value = "AssetID"
$ws.Range("D49").Value = "INTEGER"
$ws.Range("E49").Value = "SceneAssetTable"

# ---------------------------------------------------------------
# 4) Apply the grey "table name" header style (same style already
#    used by B5/B10/B14/...) to the newly introduced table headers
# ---------------------------------------------------------------
$ws.Range("B5").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B42").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B45").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B48").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B24").Value = "LightAsset"
$ws.Range("B29").Value = "MarkerAsset"
$ws.Range("B33").Value = "Scene"
$ws.Range("B42").Value = "LightAnchor"
$ws.Range("B45").Value = "MarkerAnchor"
$ws.Range("B48").Value = "SceneAssetAnchor"

# ---------------------------------------------------------------
# 5) Create + apply the new orange "Primary Key / Class Table
#    Inheritance" marker style
# ---------------------------------------------------------------
$ws.Range("C36").Interior.Color = 49407
$ws.Range("C36").Copy()
$ws.Range("F36").PasteSpecial(-4122)
$ws.Range("C36").Copy()
$ws.Range("C42").PasteSpecial(-4122)
$ws.Range("C36").Copy()
$ws.Range("C45").PasteSpecial(-4122)
$ws.Range("C36").Copy()
$ws.Range("C48").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C36").Value = "AnchorID"
$ws.Range("F36").Value = "*Shared Primary Key and Class Table Inheritance"
$ws.Range("C42").Value = "AnchorID"
$ws.Range("C45").Value = "AnchorID"
$ws.Range("C48").Value = "AnchorID"

# ---------------------------------------------------------------
# 6) Column width + selection/view state
# ---------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 17.7
$ws.Range("J25").Select()
